# Refresh the cryptos price/volume snapshot (GitHub Actions cron update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold plain text like "27.197.16" or "1.816.72" which
# Excel would otherwise auto-coerce into numbers; force text format first so
# the values round-trip exactly as strings, matching the source data feed.
$priceCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.197.16"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.851.76"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D5").Value = "313.52"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4647"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "0.3714"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "0.07284"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "0.8872"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "20.04"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "0.07829"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "1.816.72"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "5.378"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "6.522"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "90.81"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "0.000008925"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "27.227.36"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "5.074"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "10.51"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "2.135.28"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "1.954"
$ws.Range("E25").Value = "  +5.76%  "
$ws.Range("D26").Value = "151.79"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "2.036"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "115.75"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "5.055"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "0.08805"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "3.136"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "0.7672"
$ws.Range("E33").Value = "  +5.33%  "
$ws.Range("E34").Value = "  +3.32%  "
$ws.Range("D35").Value = "4.505"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").Value = "2.722"
$ws.Range("E36").Value = "  +10.63%  "
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "0.05199"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "2.936"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "7.015"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "0.5102"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "8.436"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "0.4793"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "102.84"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "1.638"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "0.06206"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  +1.20%  "
